$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D formulas: 1024/Cx -> 8*1024/Cx (D3 is a standalone formula, the
# rest share formula group si="1"; setting them this way keeps that shape
# and refreshes the cached <v> results).
$ws.Range("D3").Formula = "=8*1024/C3"
$ws.Range("D4:D11").Formula = "=8*1024/C4"

# Move the active selection from A10 to J29.
$ws.Range("J29").Select() | Out-Null
